$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing header style (from F1, which already carries the bold/border/center style)
# into the two new header cells G1 and H1, then set their text.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F1").Value = "injured_or_killed_per_10000_per_5_years"
$ws.Range("G1").Value = "injured_per_10000"
$ws.Range("H1").Value = "killed_per_10000"

# Updated population (E) and recomputed rate columns (F, G, H) for each state row.
$data = @(
    @(2, 4858319, 10.04050989653005, 6.170858685895266, 3.86965121063479),
    @(3, 737857, 8.023234854450118, 4.404647513000486, 3.61858734144963),
    @(4, 6889326, 3.178830556138583, 1.590866798871181, 1.587963757267402),
    @(5, 2984306, 7.103829164971688, 4.513612210007955, 2.590216954963733),
    @(6, 38960501, 3.389586802284704, 1.961987090463749, 1.427599711820954),
    @(7, 5485357, 3.516635289189017, 2.065499109720661, 1.451136179468355),
    @(8, 3583240, 4.462441812437905, 3.510789118228196, 0.9516526942097097),
    @(9, 944756, 11.32567562418233, 9.028786268623856, 2.296889355558472),
    @(10, 678429, 27.62263995200677, 20.85700935543734, 6.765630596569427),
    @(11, 20403663, 5.381876773793019, 3.466044307828453, 1.915832465964567),
    @(12, 10240259, 6.359214156595063, 3.960837318665475, 2.398376837929587),
    @(13, 1419563, 1.042574369718005, 0.5987758204461514, 0.4437985492718534),
    @(14, 1674001, 1.88171930602192, 0.8900831003087811, 0.9916362057131387),
    @(15, 12826826, 13.19344317916217, 10.53573191060672, 2.657711268555448),
    @(16, 6626443, 6.88755641601384, 4.460915154631225, 2.426641261382615),
    @(17, 3124141, 3.501762564493728, 2.627922363299223, 0.8738402011945043),
    @(18, 2905603, 5.069515690891013, 3.204154180732881, 1.865361510158132),
    @(19, 4432764, 6.677549267229205, 4.272729159504093, 2.404820107725112),
    @(20, 4656925, 14.12305330234006, 9.444000064420193, 4.679053237919872),
    @(21, 1331961, 1.831885468118061, 0.9910200073425572, 0.8408654607755032),
    @(22, 5988171, 8.112660777389289, 5.288760123917637, 2.823900653471653),
    @(23, 6806028, 3.19275794927673, 2.499255072121361, 0.6935028771553688),
    @(24, 9947064, 4.602363069142815, 3.00591209627283, 1.596450972869984),
    @(25, 5506968, 2.500468497365519, 1.663347235720273, 0.8371212616452466),
    @(26, 2987519, 10.23926542391864, 6.302888784975091, 3.936376638943551),
    @(27, 6080679, 9.40848875594321, 5.895723158548577, 3.512765597394633),
    @(28, 1036653, 2.797464532490621, 1.311914401443877, 1.485550131046744),
    @(29, 1897175, 4.738624533846377, 3.552650651626761, 1.185973882219616),
    @(30, 2895891, 6.043045128425068, 3.283963381218423, 2.759081747206645),
    @(31, 1340145, 1.731155957004652, 1.074510594002888, 0.6566453630017647),
    @(32, 8871970, 4.179455070294422, 2.862949266059286, 1.316505804235136),
    @(33, 2091214, 4.963624000221881, 2.596577872948441, 2.36704612727344),
    @(34, 19613918, 3.448571570453185, 2.579290889255273, 0.8692806811979127),
    @(35, 10102084, 6.784738673723164, 4.582222836396926, 2.202515837326239),
    @(36, 746826, 2.664609962695461, 1.7406999756302, 0.9239099870652602),
    @(37, 11627878, 7.061477597202172, 4.904592222243817, 2.156885374958355),
    @(38, 3906463, 6.376612296084718, 3.967783644693422, 2.408828651391297),
    @(39, 4052831, 2.812848598917645, 1.712383269867409, 1.100465329050237),
    @(40, 12786714, 5.827142141444628, 3.954104236631866, 1.873037904812761),
    @(41, 1056302, 3.871998727636604, 3.275578385726809, 0.5964203419097948),
    @(42, 4923837, 9.533215660875857, 6.263407988526021, 3.269807672349836),
    @(43, 859999, 2.372095781506723, 1.325582936724345, 1.046512844782378),
    @(44, 6625528, 9.511694766062417, 6.758706626852984, 2.752988139209434),
    @(45, 27625500, 4.036850011764493, 2.210276737072632, 1.826573274691861),
    @(46, 3018802, 2.265799479396132, 1.331654079995972, 0.9341453994001594),
    @(47, 624833, 2.080555924543038, 1.168312173012629, 0.9122437515304089),
    @(48, 8383367, 5.99401171390922, 4.253660850109508, 1.740350863799712),
    @(49, 7237383, 2.956869907257913, 1.778267089084549, 1.178602818173365),
    @(50, 1832961, 5.368362992993304, 3.540719087858389, 1.827643905134916),
    @(51, 5769906, 5.162995722980583, 3.669037242547799, 1.493958480432783),
    @(52, 581835, 2.148375398523636, 0.8937241657858327, 1.254651232737804)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]
}
